$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "現金" (cash) — add owner/total/property_category/... columns
# ---------------------------------------------------------------------
$wsCash = $wb.Worksheets.Item("現金")

# Header row (row 1) — replicate the bold/centered header style (s=1)
# from an existing header cell, then fill in the new header labels.
$wsCash.Range("B1").Copy() | Out-Null
$wsCash.Range("C1:K1").PasteSpecial(-4122) | Out-Null

$wsCash.Range("B1").Value = "currency"
$wsCash.Range("C1").Value = "owner"
$wsCash.Range("D1").Value = "total"
$wsCash.Range("E1").Value = "property_category"
$wsCash.Range("F1").Value = "category"
$wsCash.Range("G1").Value = "date"
$wsCash.Range("H1").Value = "legislator_name"
$wsCash.Range("I1").Value = "legislator_id"
$wsCash.Range("J1").Value = "source_file"
$wsCash.Range("K1").Value = "index"

# Data row (row 2) — replicate the data-row style (s=2) from an existing
# data cell, then fill in the values.
$wsCash.Range("B2").Copy() | Out-Null
$wsCash.Range("E2:K2").PasteSpecial(-4122) | Out-Null

$wsCash.Range("B2").Value = "新臺幣"
$wsCash.Range("E2").Value = "cash"
$wsCash.Range("F2").Value = "normal"
# Force text format so the date-like string is not auto-converted to a
# date serial number, then write it and restore the shared "no format"
# look (matching the other plain data cells on this row).
$wsCash.Range("G2").NumberFormat = "@"
$wsCash.Range("G2").Value = "2013-07-11"
$wsCash.Range("B2").Copy() | Out-Null
$wsCash.Range("G2").PasteSpecial(-4122) | Out-Null
$wsCash.Range("H2").Value = "林正二"
$wsCash.Range("I2").Value = 788
$wsCash.Range("J2").Value = "tmp685a1"
$wsCash.Range("K2").Value = 51

# ---------------------------------------------------------------------
# Sheet "存款" (deposit) — add owner/total/property_category/... columns
# ---------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item("存款")

# Header row (row 1)
$wsDeposit.Range("B1").Copy() | Out-Null
$wsDeposit.Range("E1:M1").PasteSpecial(-4122) | Out-Null

$wsDeposit.Range("B1").Value = "bank"
$wsDeposit.Range("C1").Value = "deposit_type"
$wsDeposit.Range("D1").Value = "currency"
$wsDeposit.Range("E1").Value = "owner"
$wsDeposit.Range("F1").Value = "total"
$wsDeposit.Range("G1").Value = "property_category"
$wsDeposit.Range("H1").Value = "category"
$wsDeposit.Range("I1").Value = "date"
$wsDeposit.Range("J1").Value = "legislator_name"
$wsDeposit.Range("K1").Value = "legislator_id"
$wsDeposit.Range("L1").Value = "source_file"
$wsDeposit.Range("M1").Value = "index"

# Data row (row 2)
$wsDeposit.Range("B2").Copy() | Out-Null
$wsDeposit.Range("G2:M2").PasteSpecial(-4122) | Out-Null

$wsDeposit.Range("B2").Value = "臺灣銀行群賢分行"
$wsDeposit.Range("C2").Value = "活期儲蓄存款"
$wsDeposit.Range("G2").Value = "deposit"
$wsDeposit.Range("H2").Value = "normal"
$wsDeposit.Range("I2").NumberFormat = "@"
$wsDeposit.Range("I2").Value = "2013-07-11"
$wsDeposit.Range("B2").Copy() | Out-Null
$wsDeposit.Range("I2").PasteSpecial(-4122) | Out-Null
$wsDeposit.Range("J2").Value = "林正二"
$wsDeposit.Range("K2").Value = 788
$wsDeposit.Range("L2").Value = "tmp685a1"
$wsDeposit.Range("M2").Value = 55
